$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Cell-by-cell replacements, ordered so that a cell whose new value equals
# another cell's old (pre-edit) value is processed first; this avoids a later
# ReplaceAll search from matching text that was only just inserted by an earlier step.
$t.Cell(1,1).Range.Find.Execute("83÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "94÷8=", 2) | Out-Null
$t.Cell(1,2).Range.Find.Execute("14÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "53÷5=", 2) | Out-Null
$t.Cell(1,3).Range.Find.Execute("48÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "11÷8=", 2) | Out-Null
$t.Cell(1,4).Range.Find.Execute("38÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "56÷3=", 2) | Out-Null
$t.Cell(1,5).Range.Find.Execute("39÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "50÷2=", 2) | Out-Null
$t.Cell(5,1).Range.Find.Execute("78÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷9=", 2) | Out-Null
$t.Cell(5,2).Range.Find.Execute("64÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷7=", 2) | Out-Null
$t.Cell(5,3).Range.Find.Execute("65÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2) | Out-Null
$t.Cell(5,4).Range.Find.Execute("35÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 2) | Out-Null
$t.Cell(9,2).Range.Find.Execute("60÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "18÷4=", 2) | Out-Null
$t.Cell(5,5).Range.Find.Execute("93÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 2) | Out-Null
$t.Cell(9,1).Range.Find.Execute("28÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "17÷8=", 2) | Out-Null
$t.Cell(9,3).Range.Find.Execute("35÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷8=", 2) | Out-Null
$t.Cell(9,4).Range.Find.Execute("90÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "27÷6=", 2) | Out-Null
$t.Cell(9,5).Range.Find.Execute("77÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷7=", 2) | Out-Null
$t.Cell(13,1).Range.Find.Execute("88÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷6=", 2) | Out-Null
$t.Cell(13,2).Range.Find.Execute("32÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "19÷9=", 2) | Out-Null
$t.Cell(13,3).Range.Find.Execute("65÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷7=", 2) | Out-Null
$t.Cell(13,4).Range.Find.Execute("96÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "98÷3=", 2) | Out-Null
$t.Cell(13,5).Range.Find.Execute("25÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "87÷6=", 2) | Out-Null
$t.Cell(17,1).Range.Find.Execute("33÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "36÷3=", 2) | Out-Null
$t.Cell(17,2).Range.Find.Execute("30÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "57÷3=", 2) | Out-Null
$t.Cell(17,3).Range.Find.Execute("83÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "79÷2=", 2) | Out-Null
$t.Cell(17,4).Range.Find.Execute("19÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷7=", 2) | Out-Null
$t.Cell(17,5).Range.Find.Execute("86÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "62÷8=", 2) | Out-Null
